$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "4x4 Squat Racks": the first product row ("The Corporate Rack" / Oak
# Club MFG) was dropped from the feed. Remove data row 2 entirely and let
# every row below it shift up by one (row 6 -> row 5, etc.).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
$ws1.Rows("2:2").Delete()

# Row deletion shifts the surviving hyperlinks' anchors correctly, but this
# host leaves the two hyperlink entries that used to point at the
# now-removed last row (old row 6) dangling instead of dropping them. Rebuild
# the hyperlinks collection from scratch with exactly the 8 links that should
# remain (rows 2-5, columns E/F), matching the shifted data.
$ws1.Hyperlinks.Delete()

$hyperlinkData = @(
  @("E2", "https://titan.fitness/cdn/shop/files/401223_01.jpg?v=1722443777&width=1946"),
  @("F2", "https://titan.fitness/products/titan-series-power-rack-90-36?variant=47930285916437"),
  @("E3", "https://garagegymlab.com/wp-content/uploads/Rogue-RM-3-Monster-Rack-2.0-Blue.jpg"),
  @("F3", "https://www.roguefitness.com/rm-3-bolt-together-monster-rack-2-0"),
  @("E4", "https://shop.straydogstrength.com/cdn/shop/files/2120-v2-FRAME.jpg?v=1739385447&width=1946"),
  @("F4", "https://shop.straydogstrength.com/products/alpha-half-rack"),
  @("E5", "https://cdn.shopify.com/s/files/1/2559/4942/products/XL_SingleRack_BlackTexture.210.jpg?v=1567697449"),
  @("F5", "https://www.sorinex.com/products/xl-half-rack?Attachment+Color=Black+Texture&Upgrades=None")
)

foreach ($entry in $hyperlinkData) {
  $ws1.Hyperlinks.Add($ws1.Range($entry[0]), $entry[1]) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Squat Stands": price refresh for "The Associate Squat Stand".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").Value = "$1,544.00"
